# Daily attendance processing - 2025-12-01 19:07:25
# For each row's "Recorded By" (column G) value that contains multiple
# comma-separated names/emails, rotate the list so the last entry moves
# to the front (others shift right by one position).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = 7
    $val = $cell.Value2

    if ($null -eq $val) {
        continue
    }

    $text = [string]$val
    if ($text -eq "") {
        continue
    }

    $parts = $text -split ", "

    if ($parts.Count -ge 2) {
        $rotated = @($parts[$parts.Count - 1]) + $parts[0..($parts.Count - 2)]
        $cell.Value = [string]::Join(", ", $rotated)
    }
}
